# "Generate Report for Handoff"
#
# b.md has been handed off again: its zh-cn / de-de status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", a fresh
# handoff file (b.63290e5768f688058c7b37413b0a5c26c308f864.<locale>.xlf)
# and handoff datetime are recorded, and the Overview sheet's rollup date
# for b.md is updated to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 is b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-27-21 02:27:21"

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 is b.md
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-21 02:27:17"

# The engine's Hyperlinks.Add always appends (never replaces in place) and
# Hyperlink.Delete on a single item is a no-op, so the only clean way to
# retarget the D3 hyperlink's display text is to clear the whole
# collection for the sheet and rebuild it in the original order, changing
# only the one entry that needs a new display string.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/a.md", [Type]::Missing, [Type]::Missing, ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e089fe9d48f3e213d316fda9a5919c65c531736/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/95859e958b26121f55b9044b012db9ca5fa57099/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/67da25b5c08014ce58e8ad221a8d4767df48d425/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/b.md", [Type]::Missing, [Type]::Missing, ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e089fe9d48f3e213d316fda9a5919c65c531736/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/95859e958b26121f55b9044b012db9ca5fa57099/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/67da25b5c08014ce58e8ad221a8d4767df48d425/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet - row 3 is b.md
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-21 02:27:21"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$dede.Hyperlinks.Add($dede.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/a.md", [Type]::Missing, [Type]::Missing, ".md")
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f43a641668f0430fd1b0a8146d4641ce6feb246/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/dd807a51527dfc103e5ac7dd6a336f8f66663ed6/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bf8497415cd106261cdcfa59cf7733af2f6ee37d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md")
$dede.Hyperlinks.Add($dede.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/b.md", [Type]::Missing, [Type]::Missing, ".md")
$dede.Hyperlinks.Add($dede.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f43a641668f0430fd1b0a8146d4641ce6feb246/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/dd807a51527dfc103e5ac7dd6a336f8f66663ed6/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bf8497415cd106261cdcfa59cf7733af2f6ee37d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
